$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing 10 rows (row 11 .. row 20)
# Columns: A=id (number), B=data (text date), C=questao_id, D=disciplina, E=assunto, F=anotacao
$rows = @(
    @{ row=11; id=10; data="27/12/2025 05:40"; questao="937"; questaoIsText=$false; disciplina="Inglês"; assunto="Pronouns"; anotacao="Retirar os números de linhas no meio do testo" },
    @{ row=12; id=11; data="27/12/2025 05:43"; questao="935"; questaoIsText=$false; disciplina="Inglês"; assunto="Semantic"; anotacao="Salvar a palavra yield no sistema de english learning. nesse caso foi um verbo que teve o mesmo significado de produce" },
    @{ row=13; id=12; data="27/12/2025 05:46"; questao="934"; questaoIsText=$false; disciplina="Inglês"; assunto="Semantic"; anotacao="Incluir palavra Hence no english learnig" },
    @{ row=14; id=13; data="27/12/2025 05:48"; questao="933"; questaoIsText=$false; disciplina="Inglês"; assunto="Interpretação de Texto"; anotacao="Analisar esse erro, marquei letra C" },
    @{ row=15; id=14; data="27/12/2025 06:28"; questao="263"; questaoIsText=$true;  disciplina="Português"; assunto="Sinônimo E Antônimo"; anotacao="Criar um flashcard do significado da palavra resignação" },
    @{ row=16; id=15; data="27/12/2025 06:35"; questao="125"; questaoIsText=$true;  disciplina="Português"; assunto="Correlação Verbal"; anotacao="Estudar esssa tal uma variação `"estranha/formal`" `"Chovesse`" com o mesmo sentido de `"choveria`"" },
    @{ row=17; id=16; data="27/12/2025 06:38"; questao="62";  questaoIsText=$true;  disciplina="Português"; assunto="Conjunção"; anotacao="Investigar o valor semântico de concessão e conformidade, marquei conformidade, letra D, nessa questão" },
    @{ row=18; id=17; data="27/12/2025 08:40"; questao="279"; questaoIsText=$true;  disciplina="Português"; assunto="Coesão"; anotacao="colocar em negrito a palavra `"que`" nessa questão" },
    @{ row=19; id=18; data="27/12/2025 08:41"; questao="262"; questaoIsText=$true;  disciplina="Português"; assunto="Sinônimo E Antônimo"; anotacao="repetida?" },
    @{ row=20; id=19; data="27/12/2025 08:43"; questao="261"; questaoIsText=$true;  disciplina="Português"; assunto="Sinônimo E Antônimo"; anotacao="Criar um flashcard de significado da palavra Reminiscência" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.row, 1).Value = $r.id
    $ws.Cells.Item($r.row, 2).Value = $r.data

    if ($r.questaoIsText) {
        # Values that were stored as text in the source data even though
        # they look numeric (e.g. pasted/exported as strings) - force text
        # the same way Excel does for an apostrophe-prefixed entry.
        $ws.Cells.Item($r.row, 3).Value = "'" + $r.questao
    } else {
        $ws.Cells.Item($r.row, 3).Value = [double]$r.questao
    }

    $ws.Cells.Item($r.row, 4).Value = $r.disciplina
    $ws.Cells.Item($r.row, 5).Value = $r.assunto
    $ws.Cells.Item($r.row, 6).Value = $r.anotacao
}
